$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2101.1428
$ws.Range("J17").Value = 2101.1428
$ws.Range("L17").Value = 6303.428400000001
$ws.Range("N17").Value = -6639.428400000001
$ws.Range("H58").Value = 3430.3125
$ws.Range("I58").Value = 393.33334
$ws.Range("K58").Value = 1180.00002
$ws.Range("M58").Value = -1030.00002
$ws.Range("H74").Value = 10651.866
$ws.Range("I74").Value = 8829.076999999999
$ws.Range("K74").Value = 8829.076999999999
$ws.Range("M74").Value = -7893.076999999999
$ws.Range("H77").Value = 10651.866
$ws.Range("I77").Value = 8829.076999999999
$ws.Range("K77").Value = 44145.38499999999
$ws.Range("M77").Value = -39465.38499999999
$ws.Range("H80").Value = 724.62067
$ws.Range("I80").Value = 740.2222
$ws.Range("J80").Value = 699.0909
$ws.Range("K80").Value = 2220.6666
$ws.Range("L80").Value = 2097.2727
$ws.Range("M80").Value = -1222.6666
$ws.Range("N80").Value = -4093.2727
$ws.Range("H83").Value = 724.62067
$ws.Range("I83").Value = 740.2222
$ws.Range("J83").Value = 699.0909
$ws.Range("K83").Value = 6661.999800000001
$ws.Range("L83").Value = 6291.8181
$ws.Range("M83").Value = -1669.999800000001
$ws.Range("N83").Value = -16275.8181
$ws.Range("H100").Value = 3333.1667
$ws.Range("I100").Value = 3500
$ws.Range("K100").Value = 3500
$ws.Range("M100").Value = -2959
$ws.Range("H111").Value = 107307.6
$ws.Range("I111").Value = 173175.83
$ws.Range("J111").Value = 8505.25
$ws.Range("K111").Value = 519527.49
$ws.Range("L111").Value = 25515.75
$ws.Range("M111").Value = -516460.49
$ws.Range("N111").Value = -31649.75
$ws.Range("H113").Value = 4477
$ws.Range("I113").Value = 3679.4
$ws.Range("K113").Value = 3679.4
$ws.Range("M113").Value = -425.4000000000001
$ws.Range("H137").Value = 2961.8572
$ws.Range("I137").Value = 2814.5334
$ws.Range("K137").Value = 8443.600199999999
$ws.Range("M137").Value = -5893.600199999999
$ws.Range("H138").Value = 4495.9614
$ws.Range("J138").Value = 5557.324
$ws.Range("L138").Value = 16671.972
$ws.Range("N138").Value = -26951.972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5064.8086
$ws.Range("I32").Value = 5064.8086
$ws.Range("K32").Value = 5064.8086
$ws.Range("M32").Value = -4777.8086
$ws.Range("H61").Value = 2114.963
$ws.Range("I61").Value = 1976.2
$ws.Range("J61").Value = 3849.5
$ws.Range("K61").Value = 1976.2
$ws.Range("L61").Value = 3849.5
$ws.Range("M61").Value = -1764.2
$ws.Range("N61").Value = -4273.5
$ws.Range("H136").Value = 2114.963
$ws.Range("I136").Value = 1976.2
$ws.Range("J136").Value = 3849.5
$ws.Range("K136").Value = 5928.6
$ws.Range("L136").Value = 11548.5
$ws.Range("M136").Value = -3378.6
$ws.Range("N136").Value = -16648.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1462
$ws.Range("I99").Value = 1327.5
$ws.Range("K99").Value = 1327.5
$ws.Range("M99").Value = 170.5
$ws.Range("H107").Value = 1690.1875
$ws.Range("I107").Value = 1275.1111
$ws.Range("J107").Value = 2223.8572
$ws.Range("K107").Value = 1275.1111
$ws.Range("L107").Value = 2223.8572
$ws.Range("M107").Value = 644.8888999999999
$ws.Range("N107").Value = -6063.8572
$ws.Range("H134").Value = 505100
$ws.Range("I134").Value = 10200
$ws.Range("K134").Value = 30600
$ws.Range("M134").Value = -28065

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 59690.61
$ws.Range("I31").Value = 2898
$ws.Range("J31").Value = 75917.07000000001
$ws.Range("K31").Value = 2898
$ws.Range("L31").Value = 75917.07000000001
$ws.Range("M31").Value = -2603
$ws.Range("N31").Value = -76507.07000000001
$ws.Range("H34").Value = 59690.61
$ws.Range("I34").Value = 2898
$ws.Range("J34").Value = 75917.07000000001
$ws.Range("K34").Value = 2898
$ws.Range("L34").Value = 75917.07000000001
$ws.Range("M34").Value = -2696
$ws.Range("N34").Value = -76321.07000000001
$ws.Range("H58").Value = 4385.32
$ws.Range("I58").Value = 3543.9211
$ws.Range("J58").Value = 7049.75
$ws.Range("K58").Value = 3543.9211
$ws.Range("L58").Value = 7049.75
$ws.Range("M58").Value = -3340.9211
$ws.Range("N58").Value = -7455.75
$ws.Range("H76").Value = 10000
$ws.Range("I76").Value = 10000
$ws.Range("K76").Value = 10000
$ws.Range("M76").Value = -9685
$ws.Range("H79").Value = 10000
$ws.Range("I79").Value = 10000
$ws.Range("K79").Value = 10000
$ws.Range("M79").Value = -8908
$ws.Range("H94").Value = 3633.3333
$ws.Range("J94").Value = 3633.3333
$ws.Range("L94").Value = 3633.3333
$ws.Range("N94").Value = -4535.3333
$ws.Range("H105").Value = 561.2353000000001
$ws.Range("I105").Value = 582.25
$ws.Range("K105").Value = 582.25
$ws.Range("M105").Value = 1164.75
$ws.Range("H107").Value = 578
$ws.Range("I107").Value = 301.2857
$ws.Range("J107").Value = 1062.25
$ws.Range("K107").Value = 301.2857
$ws.Range("L107").Value = 1062.25
$ws.Range("M107").Value = 1618.7143
$ws.Range("N107").Value = -4902.25
$ws.Range("H132").Value = 1143.7858
$ws.Range("I132").Value = 963.5454999999999
$ws.Range("K132").Value = 2890.6365
$ws.Range("M132").Value = -360.6364999999996
$ws.Range("H136").Value = 4385.32
$ws.Range("I136").Value = 3543.9211
$ws.Range("J136").Value = 7049.75
$ws.Range("K136").Value = 10631.7633
$ws.Range("L136").Value = 21149.25
$ws.Range("M136").Value = -8081.763300000001
$ws.Range("N136").Value = -26249.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2511647.5
$ws.Range("I4").Value = 2733839.2
$ws.Range("J4").Value = 845210
$ws.Range("K4").Value = 8201517.600000001
$ws.Range("L4").Value = 2535630
$ws.Range("M4").Value = -8201405.600000001
$ws.Range("N4").Value = -2535854
$ws.Range("H13").Value = 3018.5
$ws.Range("J13").Value = 4212.5
$ws.Range("L13").Value = 12637.5
$ws.Range("N13").Value = -12973.5
$ws.Range("H141").Value = 3166.3333
$ws.Range("I141").Value = 3166.3333
$ws.Range("K141").Value = 9498.999899999999
$ws.Range("M141").Value = -4318.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 12450
$ws.Range("J44").Value = 12450
$ws.Range("L44").Value = 12450
$ws.Range("N44").Value = -13642
$ws.Range("H113").Value = 483988.06
$ws.Range("J113").Value = 10739.286
$ws.Range("L113").Value = 10739.286
$ws.Range("N113").Value = -15079.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1335.1305
$ws.Range("I22").Value = 657.1429000000001
$ws.Range("K22").Value = 657.1429000000001
$ws.Range("M22").Value = -362.1429000000001
$ws.Range("H27").Value = 1335.1305
$ws.Range("I27").Value = 657.1429000000001
$ws.Range("K27").Value = 657.1429000000001
$ws.Range("M27").Value = -550.1429000000001
$ws.Range("H40").Value = 171667.17
$ws.Range("I40").Value = 335334.34
$ws.Range("K40").Value = 335334.34
$ws.Range("M40").Value = -335198.34
$ws.Range("H55").Value = 879.4167
$ws.Range("I55").Value = 221.86667
$ws.Range("K55").Value = 221.86667
$ws.Range("M55").Value = -48.86667
$ws.Range("H93").Value = 4569
$ws.Range("I93").Value = 4499
$ws.Range("K93").Value = 4499
$ws.Range("M93").Value = -3251

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
